# Actualización automática 2025-06-01 08:00:06
#
# This script applies the recorded edit to the workbook:
#  1) Sheet "VENTAS POR GRUPO":
#     - Column I (9) width changes from 12 to 9 characters.
#     - All the monthly figures in C2:N263 are reset to 0.
#     - The summary row 264 ("N de 262" counts) is updated to "0 de 262"
#       for every column, matching the now-all-zero data above.
#  2) Sheet "VENTA MENSUAL":
#     - Column F (6) width changes from 14 to 11 characters.
#     - The month headers roll forward by one month
#       (febrero/marzo/abril/mayo -> marzo/abril/mayo/junio).
#     - Every data row's figures roll forward one column to the left
#       (C<-D, D<-E, E<-F) with the new trailing column F filled with 0,
#       matching the month header shift.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Column I width: 12 -> 9 (Excel re-adds ~0.8333 chars of padding when the
# ColumnWidth property is read back, so subtract that offset on the way in).
$ws1.Columns.Item(9).ColumnWidth = 9 - 5/6

# Zero out every monthly value in the data block (rows 2-263, columns C-N).
$ws1.Range("C2:N263").Value = 0

# Update the "N de 262" summary row to reflect that every count is now 0.
$ws1.Range("C264:N264").Value = "0 de 262"

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column F width: 14 -> 11 (same padding offset as above).
$ws2.Columns.Item(6).ColumnWidth = 11 - 5/6

# Shift the month header labels one column to the left, introducing "junio".
$ws2.Range("C1").Value = "marzo"
$ws2.Range("D1").Value = "abril"
$ws2.Range("E1").Value = "mayo"
$ws2.Range("F1").Value = "junio"

# Shift every data row (2-264, including the totals row) one column to the
# left, filling the newly vacated column F with 0.
$dataRange = $ws2.Range("C2:F264")
$vals = $dataRange.Value2
$nrows = $vals.GetUpperBound(0)
$ncols = $vals.GetUpperBound(1)
for ($r = 1; $r -le $nrows; $r++) {
    for ($c = 1; $c -lt $ncols; $c++) {
        $vals[$r, $c] = $vals[$r, $c + 1]
    }
    $vals[$r, $ncols] = 0
}
$dataRange.Value2 = $vals
